# Applies the daily crypto price / volume(1h) refresh described in the
# commit diff ("Updated cryptos list ... with GitHub Actions").
#
# Note on quoting: several Price values (column D) are plain decimals
# (e.g. "0.999", "584.12"). The source workbook stores every Price/
# Volume cell as text (inline string), so a leading apostrophe is used
# for those values to force Excel to keep them as text instead of
# auto-converting them to numbers -- exactly what typing them into
# Excel with an existing text-looking column would do. Values that
# already contain extra separators (e.g. "62.406.89") or non-numeric
# characters are never auto-converted, so they are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-38: refresh Price (D) and Volume(1h) (E) ---
$ws.Range('D2').Value = '62.406.89'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '3.013.43'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''584.12'
$ws.Range('E5').Value = '  -1.60%  '
$ws.Range('D6').Value = '''145.88'
$ws.Range('E6').Value = '  -5.97%  '
$ws.Range('D8').Value = '''0.527'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '3.009.99'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('D10').Value = '''0.149'
$ws.Range('E10').Value = '  -4.88%  '
$ws.Range('D11').Value = '''5.78'
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('D12').Value = '''0.462'
$ws.Range('D13').Value = '''0.0000229'
$ws.Range('E13').Value = '  -3.84%  '
$ws.Range('D14').Value = '''34.53'
$ws.Range('E14').Value = '  -6.63%  '
$ws.Range('D16').Value = '3.500.12'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '''7.13'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '62.357.12'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = '3.007.62'
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('D20').Value = '''456.91'
$ws.Range('E20').Value = '  -7.20%  '
$ws.Range('D21').Value = '''13.97'
$ws.Range('D22').Value = '''0.689'
$ws.Range('E22').Value = '  -2.72%  '
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '''81.64'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').Value = '''12.41'
$ws.Range('E25').Value = '  -3.71%  '
$ws.Range('D26').Value = '''2.23'
$ws.Range('E26').Value = '  -9.40%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '''10.05'
$ws.Range('E28').Value = '  -6.75%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -3.07%  '
$ws.Range('D31').Value = '''7.02'
$ws.Range('E31').Value = '  -5.29%  '
$ws.Range('E32').Value = '  -5.98%  '
$ws.Range('D33').Value = '''28.09'
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('D35').Value = '0.0₃0805'
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').Value = '''1.03'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('D37').Value = '''5.77'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('E38').Value = '  -5.06%  '

# --- Rows 39/40: OKB and Cosmos swap ranking positions, with refreshed values ---
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '''50.25'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '''9.16'
$ws.Range('E40').Value = '  -1.19%  '

# --- Rows 41-51: refresh Price (D) and Volume(1h) (E) ---
$ws.Range('D41').Value = '''2.91'
$ws.Range('E41').Value = '  -13.22%  '
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('D43').Value = '''391.32'
$ws.Range('E43').Value = '  -11.00%  '
$ws.Range('D44').Value = '''0.0358'
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('D45').Value = '''0.269'
$ws.Range('E45').Value = '  -7.87%  '
$ws.Range('D46').Value = '2.727.83'
$ws.Range('E46').Value = '  -4.10%  '
$ws.Range('D47').Value = '''37.28'
$ws.Range('E47').Value = '  -5.03%  '
$ws.Range('D48').Value = '''128.89'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('E51').Value = '  -1.92%  '

